# Remove stale footnote markers (e.g. " [1]", " [5]") from vaccine/category
# names and collapse the embedded line-breaks used in several BrandName /
# Vaccine labels into single spaces, per the "done !" cleanup commit.
#
# Work directly against cell text so the shared-strings table is rebuilt
# (and de-duplicated / compacted) by the engine on save - this naturally
# drops the now-unused "Afluria`nQuadrivalent" entry and folds its two
# former references onto the existing "Afluria Quadrivalent" string.

$wb = $excel.ActiveWorkbook

$replacements = [ordered]@{
    "Dengue Tetravalent Vaccine, Live [7]" = "Dengue Tetravalent Vaccine, Live "
    "DTaP [1]" = "DTaP "
    "DTaP-IPV [2]" = "DTaP-IPV "
    "DTaP-Hep B-IPV [4]" = "DTaP-Hep B-IPV "
    "DTaP-IP-HI [4]" = "DTaP-IP-HI "
    "DTaP-IPV-HIB-HEPB [6]" = "DTaP-IPV-HIB-HEPB "
    "e-IPV [5]" = "e-IPV "
    "Hepatitis A Pediatric [5]" = "Hepatitis A Pediatric "
    "Hepatitis A-Hepatitis B 18 only [3]" = "Hepatitis A-Hepatitis B 18 only "
    "Hepatitis B [5]`nPediatric/Adolescent" = "Hepatitis B  Pediatric/Adolescent"
    "Hib [5]" = "Hib "
    "HPV - Human Papillomavirus 9-valent [5]" = "HPV - Human Papillomavirus 9-valent "
    "MENB - Meningococcal Group B [5]" = "MENB - Meningococcal Group B "
    "Meningococcal Conjugate (Groups A, C, W and Y) [5]" = "Meningococcal Conjugate (Groups A, C, W and Y) "
    "Meningococcal Conjugate (Groups A, C, Y and W-135) [5]" = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
    "Measles, Mumps and Rubella (MMR) [1]" = "Measles, Mumps and Rubella (MMR) "
    "MMR/Varicella [2]" = "MMR/Varicella "
    "Pneumococcal`n13-valent [5] (Pediatric)" = "Pneumococcal 13-valent  (Pediatric)"
    "Rotavirus, Live, Oral, Pentavalent [5]" = "Rotavirus, Live, Oral, Pentavalent "
    "Rotavirus, Live, Oral, Oral [5]" = "Rotavirus, Live, Oral, Oral "
    "Tetanus and Diphtheria Toxoids [3]" = "Tetanus and Diphtheria Toxoids "
    "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis [1]" = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
    "Varicella [5]" = "Varicella "
    "Hepatitis A Adult [5]" = "Hepatitis A Adult "
    "Hepatitis A-Hepatitis B Adult [3]" = "Hepatitis A-Hepatitis B Adult "
    "Hepatitis B Adult [5]" = "Hepatitis B Adult "
    "Hepatitis B-Adult [5]" = "Hepatitis B-Adult "
    "HPV-Human Papillomavirus 9 Valent [5]" = "HPV-Human Papillomavirus 9 Valent "
    "Measles, Mumps,  Rubella [1]" = "Measles, Mumps,  Rubella "
    "Pneumococcal`n15-valent [5]" = "Pneumococcal 15-valent "
    "Pneumococcal`n20-valent [5]" = "Pneumococcal 20-valent "
    "Influenza [5]`n(Age 6 months and older)" = "Influenza  (Age 6 months and older)"
    "Fluzone`nQuadrivalent" = "Fluzone Quadrivalent"
    "Fluarix`nQuadrivalent" = "Fluarix Quadrivalent"
    "FluLaval`nQuadrivalent" = "FluLaval Quadrivalent"
    "Influenza [5]`n(Age 36 months and older)" = "Influenza  (Age 36 months and older)"
    "Influenza [5]`nLive, Intranasal (Age 2-49 years)" = "Influenza  Live, Intranasal (Age 2-49 years)"
    "FluMist`nQuadrivalent" = "FluMist Quadrivalent"
    "Afluria`nQuadrivalent" = "Afluria Quadrivalent"
}

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $ur = $ws.UsedRange
    $rows = $ur.Rows.Count
    $cols = $ur.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value()
            if ($val -ne $null -and $replacements.Contains($val)) {
                $cell.Value = $replacements[$val]
            }
        }
    }
}
